# Apply cryptos list update (prices/volumes refreshed, two coin pairs swapped rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.000.65"
$ws.Range("E2").Value = "  -2.14%  "
$ws.Range("D3").Value = "2.016.01"
$ws.Range("E3").Value = "  -3.58%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.651"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.36"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "58.82"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("E10").Value = "  -2.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0732"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.82%  "
$ws.Range("E12").Value = "  -4.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.889"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.20%  "
$ws.Range("E14").Value = "  -6.64%  "
$ws.Range("D15").Value = "2.313.11"
$ws.Range("E15").Value = "  -3.34%  "
$ws.Range("E16").Value = "  -4.93%  "
$ws.Range("D17").Value = "2.024.44"
$ws.Range("E17").Value = "  -3.54%  "
$ws.Range("D18").Value = "35.949.27"
$ws.Range("E18").Value = "  -2.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.51%  "
$ws.Range("D21").Value = "0.0₃0844"
$ws.Range("E21").Value = "  -4.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "235.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.08%  "
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("E25").Value = "  -2.40%  "
$ws.Range("E26").Value = "  +4.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.79%  "
$ws.Range("E30").Value = "  -2.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.76%  "
$ws.Range("E33").Value = "  -2.94%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0899"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.46%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.28"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -9.01%  "
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("E37").Value = "  -1.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -10.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.35%  "
$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0212"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.71%  "
$ws.Range("E43").Value = "  -5.77%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "91.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.35%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.385.72"
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0889"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.49%  "
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("E50").Value = "  -8.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.72%  "
